$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (the "2007年" data row). This shifts rows 3-6 up to become rows 2-5.
$ws.Rows.Item(2).Delete()
